# chore: adapt column header formatting to respective input file names
#
# - Rename the "_old"/"_new" header suffixes to the actual format-version
#   identifiers they represent ("_FV2404" / "_FV2410").
# - Freeze the header row so it stays visible while scrolling.
# - Turn the sheet's data range into a proper Excel Table (ListObject) so the
#   renamed headers double as the table's column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells --------------------------------------------
# Columns A-J describe the "old" (FV2404) side of the comparison, column K
# is the literal "diff" column, and columns L-U describe the "new" (FV2410)
# side.
$oldSideCols = @("A1", "B1", "C1", "D1", "E1", "F1", "G1", "H1", "I1", "J1")
foreach ($addr in $oldSideCols) {
    $cell = $ws.Range($addr)
    $current = [string]$cell.Value()
    $cell.Value = $current.Replace("_old", "_FV2404")
}

$newSideCols = @("L1", "M1", "N1", "O1", "P1", "Q1", "R1", "S1", "T1", "U1")
foreach ($addr in $newSideCols) {
    $cell = $ws.Range($addr)
    $current = [string]$cell.Value()
    $cell.Value = $current.Replace("_new", "_FV2410")
}

# --- 2. Freeze the header row (split below row 1) ----------------------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Wrap the data range in an Excel Table ---------------------------
$dataRange = $ws.UsedRange
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
